# Update the "timestamp" column (Z) for all data rows (2-112) of the
# pcsmote sample log sheet. The underlying notebook was re-run and every
# logged row got a fresh timestamp; rows that were produced within the
# same processing batch share the same timestamp value, which is why the
# update below is expressed as a handful of contiguous Z-column ranges
# rather than 111 individual cell writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z7").Value   = "2025-10-17T07:09:35.566715"
$ws.Range("Z8:Z18").Value  = "2025-10-17T07:09:35.567712"
$ws.Range("Z19:Z29").Value = "2025-10-17T07:09:35.568712"
$ws.Range("Z30:Z40").Value = "2025-10-17T07:09:35.569712"
$ws.Range("Z41:Z45").Value = "2025-10-17T07:09:35.570712"
$ws.Range("Z46:Z71").Value = "2025-10-17T07:09:35.625263"
$ws.Range("Z72:Z74").Value = "2025-10-17T07:09:35.629775"
$ws.Range("Z75:Z102").Value = "2025-10-17T07:09:35.691840"
$ws.Range("Z103:Z112").Value = "2025-10-17T07:09:35.763480"
